$wb = $excel.ActiveWorkbook

# Fix typo in the "bar" rule body on the Rules sheet: "retirn" -> "return"
$rulesSheet = $wb.Worksheets.Item("Rules")
$rulesSheet.Range("B14").Value = "return bar.foo;"

# Reflect the final active sheet/selection as left by the edit
$rulesSheet.Activate()
$rulesSheet.Range("B15").Select()
